# observed new factor values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap columns E and F (header + 13 data rows: rows 1-13)
for ($r = 1; $r -le 13; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 5).Value2 = $fVal
    $ws.Cells.Item($r, 6).Value2 = $eVal
}

# Swap the (bestFit) column widths of E and F to follow the content swap:
# before: E=16.83203125, F=24  -> after: E=24, F=16.83203125
$ws.Columns.Item(5).ColumnWidth = 23.166666666666668
$ws.Columns.Item(6).ColumnWidth = 16

# 2. Re-style rows 14:15 to match the "highlighted" rows (same style as A2:C3, yellow fill, no bold)
$srcStyle = $ws.Range("A2:C2")
$dstStyle = $ws.Range("A14:C15")
$srcStyle.Copy()
$dstStyle.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Update A14:C15 values
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "div_factor"
$ws.Range("C14").Value = 100

$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "final_div_factor"
$ws.Range("C15").Value = 100

# 4. Add D14:G15 new values
$ws.Range("D14").Value = 77.83
$ws.Range("E14").Value = 72.55
$ws.Range("F14").Value = 91.72
$ws.Range("G14").Value = 84.45

$ws.Range("D15").Value = 77.83
$ws.Range("E15").Value = 72.55
$ws.Range("F15").Value = 91.72
$ws.Range("G15").Value = 84.45

# 5. Set selection to G11
$ws.Range("G11").Select()
